# Weekly update to the fruit/vegetable price sheet: a new record was
# inserted at row 313 (pushing the existing rows 313-370 down to 314-371).
#
# Columns (A..T):
#  A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo,
#  G Producto ID, H Producto, I Categoria ID, J Categoria, K Variedad,
#  L Calidad, M Volumen, N Precio minimo, O Precio maximo,
#  P Precio promedio ponderado, Q Unidad de comercializacion,
#  R Origen, S Precio $/Kg, T Kg / unidad

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 313; Excel automatically shifts the
# former rows 313..370 down to 314..371 and carries their formatting along.
$ws.Rows.Item(313).Insert()

# Populate the new row 313 with the new record's data.
$ws.Cells.Item(313, 1).Value2 = 10
$ws.Cells.Item(313, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(313, 3).Value2 = "La Araucanía"
$ws.Cells.Item(313, 4).Value2 = 44511
$ws.Cells.Item(313, 5).Value2 = 9
$ws.Cells.Item(313, 6).Value2 = "Fruta"
$ws.Cells.Item(313, 7).Value2 = 100108
$ws.Cells.Item(313, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(313, 9).Value2 = 100108006
$ws.Cells.Item(313, 10).Value2 = "Plátano"
$ws.Cells.Item(313, 11).Value2 = "Sin especificar"
$ws.Cells.Item(313, 12).Value2 = "Pintón"
$ws.Cells.Item(313, 13).Value2 = 1755
$ws.Cells.Item(313, 14).Value2 = 18000
$ws.Cells.Item(313, 15).Value2 = 20000
$ws.Cells.Item(313, 16).Value2 = 18513
$ws.Cells.Item(313, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(313, 18).Value2 = "Ecuador"
$ws.Cells.Item(313, 19).Value2 = 926
$ws.Cells.Item(313, 20).Value2 = 20
